$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{2=15; 3=20; 4=24; 5=13; 6=17; 7=27; 8=22; 9=5; 10=2; 11=12; 12=7; 13=23; 14=11; 15=19; 16=9; 17=14; 18=16; 19=18; 20=10; 21=25; 22=6; 23=26; 24=21; 25=8; 26=4; 27=3}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot original values before overwriting anything
$orig = @{}
foreach ($r in $map.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

foreach ($dst in $map.Keys) {
    $src = $map[$dst]
    foreach ($c in $cols) {
        $ws.Range("$c$dst").Value2 = $orig[$src][$c]
    }
}
